# FAAS-459 iOS phone Device testing
#
# The test data sheet gains a new column (AL) with header
# "discDateWithTimestamp" (a new field alongside the existing
# "discDate" column), pushing the used range from A1:AK2 to A1:AL2.
# The active selection in the sheet also moves to AM6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell - this also grows the shared-string table,
# the sheet dimension and the row "spans" attributes automatically.
$ws.Range("AL1").Value = "discDateWithTimestamp"

# Move/record the current selection to AM6, matching the saved
# workbook's cursor position.
$ws.Range("AM6").Select()
